$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '28.415.17'
$ws.Range("E2").Value = '  -3.32%  '
$ws.Range("D3").Value = '1.955.96'
$ws.Range("E3").Value = '  -1.81%  '
$ws.Range("E4").Value = '  -0.64%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '321.09'
$ws.Range("E5").Value = '  -2.26%  '
$ws.Range("E6").Value = '  -0.52%  '
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.4762'
$ws.Range("E7").Value = '  -4.56%  '
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.4054'
$ws.Range("E8").Value = '  -3.69%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '53.50'
$ws.Range("E9").Value = '  -1.26%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.08467'
$ws.Range("E10").Value = '  -5.19%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '1.057'
$ws.Range("E11").Value = '  -4.83%  '
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '22.18'
$ws.Range("E12").Value = '  -4.40%  '
$ws.Range("D13").Value = '1.955.25'
$ws.Range("E13").Value = '  -4.75%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '7.603'
$ws.Range("E14").Value = '  -4.56%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '6.141'
$ws.Range("E15").Value = '  -4.53%  '
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '1.010'
$ws.Range("E16").Value = '  -0.67%  '
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '89.53'
$ws.Range("E17").Value = '  -4.25%  '
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '0.00001066'
$ws.Range("E18").Value = '  -3.65%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '0.06598'
$ws.Range("E19").Value = '  -1.66%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '18.58'
$ws.Range("E20").Value = '  -4.46%  '
$ws.Range("E21").Value = '  -0.64%  '
$ws.Range("E22").Value = '  -2.03%  '
$ws.Range("D23").Value = '28.421.14'
$ws.Range("E23").Value = '  -3.54%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '11.51'
$ws.Range("E24").Value = '  -3.71%  '
$ws.Range("E25").Value = '  -0.66%  '
$ws.Range("D26").Value = '2.189.38'
$ws.Range("E26").Value = '  -3.84%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '154.78'
$ws.Range("E27").Value = '  -1.82%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '20.20'
$ws.Range("E28").Value = '  -2.59%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '5.910'
$ws.Range("E29").Value = '  -5.37%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '2.154'
$ws.Range("E30").Value = '  -6.27%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '123.56'
$ws.Range("E31").Value = '  -3.20%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '0.9771'
$ws.Range("E32").Value = '  -7.16%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '0.09584'
$ws.Range("E33").Value = '  -3.18%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '1.440'
$ws.Range("E34").Value = '  -7.04%  '
$ws.Range("E35").Value = '  -4.01%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '5.579'
$ws.Range("E36").Value = '  -4.12%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '0.02324'
$ws.Range("E37").Value = '  -5.09%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '8.820'
$ws.Range("E38").Value = '  -4.28%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '0.06203'
$ws.Range("E39").Value = '  -2.32%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '1.249'
$ws.Range("E40").Value = '  -3.32%  '
$ws.Range("E41").Value = '  -5.35%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '11.12'
$ws.Range("E42").Value = '  -3.90%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '1.008'
$ws.Range("E43").Value = '  -0.60%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '0.1916'
$ws.Range("E44").Value = '  -5.83%  '
$ws.Range("E45").Value = '  +3.21%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '0.5954'
$ws.Range("E46").Value = '  -5.83%  '
$ws.Range("E47").Value = '  -3.45%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '2.054'
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '3.395'
$ws.Range("E49").Value = '  -2.97%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '0.00000000328'
$ws.Range("E50").Value = '  -1.29%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '0.06802'
$ws.Range("E51").Value = '  -2.17%  '
